# Edit script: updates GDP column (C) with new LSTM-predicted values,
# flips several "Colony" (AL) flags from 0 to 1, and renames several
# header labels (shared strings) to shorter/updated names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames ---
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) updated values, per row ---
$cValues = @{
    2 = 2934.187009790061
    3 = 2870.311589353206
    4 = 697.6889104500298
    5 = 1873.394108966653
    6 = 1460.056109840828
    7 = 5191.140356354663
    8 = 1909.084588129339
    9 = 6128.19547247793
    10 = 4729.735976516416
    11 = 14239.03920301361
    12 = 993.3829437244538
    13 = 4744.762791189912
    14 = 2812.435974421079
    15 = 1036.533951644687
    16 = 473.2998774917226
    17 = 2983.242707849043
    18 = 2898.942214704482
    19 = 665.6274194933962
    20 = 5555.389721901988
    21 = 1955.461557360978
    22 = 12808.034586422
    23 = 5082.354756663512
    24 = 13825.35808833117
    25 = 1357.563719132622
    26 = 5076.340174387075
    27 = 492.3430015592067
    28 = 1303.425880277445
    29 = 1037.747039954749
    30 = 982.980837581714
    31 = 8390.479071096475
    32 = 3083.80337578809
    33 = 2965.153206179127
    34 = 691.8942672110555
    35 = 1577.487171555845
    36 = 5660.517066940175
    37 = 10883.31535948899
    38 = 2024.117324382548
    39 = 6711.616186806423
    40 = 5360.226632400601
    41 = 33888.49231534224
    42 = 992.8781394745556
    43 = 5437.877690816224
    44 = 5325.160106166602
    45 = 1469.177610078392
    46 = 495.763971160512
    47 = 1057.667740311969
    48 = 683.460336640684
    49 = 692.4450379203138
    50 = 3156.723844635973
    51 = 2999.422762626143
    52 = 701.4459636783288
    53 = 1657.651524528445
    54 = 4394.543881413723
    55 = 5642.578115155247
    56 = 5745.422744292303
    57 = 5558.73713990153
    58 = 5710.587873377512
    59 = 1544.619247249133
    60 = 503.3023574516347
    61 = 1102.527430026863
    62 = 698.3833464078615
    63 = 711.0361291687414
    64 = 3212.740625904757
    65 = 3056.152683606517
    66 = 720.7128711178943
    67 = 1716.389195271215
    68 = 4699.493713911862
    69 = 5919.20956823756
    70 = 5955.175904294275
    71 = 1140.447753778042
    72 = 2286.013198234259
    73 = 1401.753174264641
    74 = 6103.590270484282
    75 = 7449.08671983612
    76 = 3008.669179463094
    77 = 6255.426161047989
    78 = 3252.634165082374
    79 = 3137.260298393558
    80 = 730.3063521039821
    81 = 558.2093442539386
    82 = 1257.483615623398
    83 = 711.3043470146426
    84 = 1775.027517189621
    85 = 4861.287098802361
    86 = 5996.49696468919
    87 = 6301.696269820412
    88 = 1128.996380299766
    89 = 2361.056581219794
    90 = 1441.783971398429
    91 = 6500.281937297324
    92 = 7580.275568826287
    93 = 3012.536723186288
    94 = 6522.736799041846
    95 = 3314.741082534716
    96 = 3210.869677115934
    97 = 729.1196658666737
    98 = 579.0880693780265
    99 = 1335.203871985052
    100 = 731.9993357350996
    101 = 1836.014008604312
    102 = 4944.191641077407
    103 = 6114.227214287786
    104 = 6661.86504232374
    105 = 1134.924536209078
    106 = 1469.192636109792
    107 = 2854.757682901436
    108 = 5176.058803160127
    109 = 3382.563653843273
    110 = 3242.636921959078
    111 = 729.8559996981501
    112 = 1875.732161108182
    113 = 1402.276549638797
    114 = 1895.214690888655
    115 = 1117.517554619445
    116 = 6262.368904654469
    117 = 22666.28622740736
}

$alRows = @(4, 7, 19, 20, 34, 36, 52, 56, 66, 70, 80, 87, 97, 104, 111)

foreach ($row in $cValues.Keys) {
    $ws.Range("C$row").Value = $cValues[$row]
}

# --- Column AL (Colony) flag flips: 0 -> 1 ---
foreach ($row in $alRows) {
    $ws.Range("AL$row").Value = 1
}
